$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.74000000000058
$ws.Range("G2").Value = 0.001886677224974487
$ws.Range("H2").Value = 0.006251404355097701
$ws.Range("K2").Value = 5.806617300348369
$ws.Range("L2").Value = "[1.729857434801854, 9.883377165894883]"
$ws.Range("M2").Value = 0.005482015466018231
$ws.Range("N2").Value = 0.005482015466018231
$ws.Range("O2").Value = -1.358526552903695
$ws.Range("P2").Value = "[-2.2013161736865428, -0.5157369321208467]"
$ws.Range("Q2").Value = 0.001722998105908813
$ws.Range("R2").Value = 0.001722998105908813
$ws.Range("S2").Value = 14.93832568517827
$ws.Range("T2").Value = "[12.675768015375517, 17.200883354981023]"
$ws.Range("W2").Value = 5.565405405405532
$ws.Range("X2").Value = 2.112792792792838
$ws.Range("Y2").Value = 9.018018018018225

# Row 3 updates
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 22.54000000000008
$ws.Range("G3").Value = [double]"6.010938985612224e-05"
$ws.Range("H3").Value = 0.0008291913948715868
$ws.Range("I3").Value = ""
$ws.Range("K3").Value = 7.044416417028114
$ws.Range("L3").Value = "[3.1908572094323784, 10.89797562462385]"
$ws.Range("M3").Value = 0.0003857109660845293
$ws.Range("N3").Value = 0.0007714219321690585
$ws.Range("O3").Value = 1.817658211986887
$ws.Range("P3").Value = "[1.1887107337907326, 2.446605690183042]"
$ws.Range("Q3").Value = [double]"3.667313719368792e-08"
$ws.Range("R3").Value = [double]"7.334627438737584e-08"
$ws.Range("S3").Value = 14.21188618180092
$ws.Range("T3").Value = "[12.052437888559396, 16.37133447504244]"
$ws.Range("W3").Value = 16.01941941941948
$ws.Range("X3").Value = 13.76316316316321
$ws.Range("Y3").Value = 18.27567567567574
